$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.405.07'
$ws.Range("E2").Value = '  -3.22%  '

$ws.Range("D3").Value = '3.520.44'
$ws.Range("E3").Value = '  -4.92%  '

$ws.Range("E4").Value = '  -0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '583.25'
$c.ClearFormats()
$ws.Range("E5").Value = '  -1.27%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '174.01'
$c.ClearFormats()
$ws.Range("E6").Value = '  -3.71%  '

$ws.Range("E7").Value = '  +0.73%  '

$ws.Range("D8").Value = '3.511.74'
$ws.Range("E8").Value = '  -4.82%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("E10").Value = '  -6.23%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.79'
$c.ClearFormats()
$ws.Range("E11").Value = '  +4.60%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.597'
$c.ClearFormats()
$ws.Range("E12").Value = '  -2.79%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '47.10'
$c.ClearFormats()
$ws.Range("E13").Value = '  -5.80%  '

$ws.Range("E14").Value = '  -3.94%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '677.81'
$c.ClearFormats()
$ws.Range("E15").Value = '  -1.12%  '

$ws.Range("D16").Value = '4.078.32'
$ws.Range("E16").Value = '  -4.89%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '8.77'
$c.ClearFormats()
$ws.Range("E17").Value = '  -3.05%  '

$ws.Range("D18").Value = '69.352.54'
$ws.Range("E18").Value = '  -3.38%  '

$ws.Range("D19").Value = '3.520.67'
$ws.Range("E19").Value = '  -4.97%  '

$ws.Range("E20").Value = '  -1.30%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.50'
$c.ClearFormats()
$ws.Range("E21").Value = '  -3.77%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '11.23'
$c.ClearFormats()
$ws.Range("E22").Value = '  -4.00%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.906'
$c.ClearFormats()
$ws.Range("E23").Value = '  -4.27%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '16.18'
$c.ClearFormats()
$ws.Range("E24").Value = '  -9.37%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '98.17'
$c.ClearFormats()
$ws.Range("E25").Value = '  -5.56%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.88'
$c.ClearFormats()
$ws.Range("E26").Value = '  -4.41%  '

$ws.Range("E27").Value = '  -0.70%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.67'
$c.ClearFormats()
$ws.Range("E29").Value = '  -6.23%  '

$ws.Range("E30").Value = '  -7.97%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '33.10'
$c.ClearFormats()
$ws.Range("E31").Value = '  -7.00%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '8.77'
$c.ClearFormats()
$ws.Range("E32").Value = '  -5.57%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.21'
$c.ClearFormats()
$ws.Range("E33").Value = '  -7.67%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '7.33'
$c.ClearFormats()
$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("E35").Value = '  -6.20%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '596.55'
$c.ClearFormats()
$ws.Range("E36").Value = '  +5.59%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.62'
$c.ClearFormats()
$ws.Range("E37").Value = '  -15.18%  '

$ws.Range("E38").Value = '  -3.50%  '

$ws.Range("E39").Value = '  -4.60%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '57.27'
$c.ClearFormats()
$ws.Range("E40").Value = '  -4.13%  '

$ws.Range("E41").Value = '  +0.13%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0441'
$c.ClearFormats()
$ws.Range("E42").Value = '  -5.62%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.337'
$c.ClearFormats()
$ws.Range("E43").Value = '  -4.53%  '

$ws.Range("E44").Value = '  -6.00%  '

$ws.Range("D45").Value = '3.420.49'
$ws.Range("E45").Value = '  -8.93%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '33.45'
$c.ClearFormats()
$ws.Range("E46").Value = '  -6.36%  '

$ws.Range("D47").Value = '0.0₃0711'
$ws.Range("E47").Value = '  -8.95%  '

$ws.Range("E48").Value = '  +0.99%  '

$ws.Range("E49").Value = '  -7.01%  '

$ws.Range("E50").Value = '  -0.61%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '5.84'
$c.ClearFormats()
$ws.Range("E51").Value = '  +19.05%  '
